$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.53"
$ws.Range("G2").Value = "'13"
$ws.Range("D3").Value = "'25.00"
$ws.Range("G3").Value = "'13"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").Value = "'5.179"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("G4").Value = "'13"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").Value = "'0.05731"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("G5").Value = "'13"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'6.470"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("G6").Value = "'13"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.066"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("G7").Value = "'13"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8095"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("G8").Value = "'13"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8405"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("G9").Value = "'13"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1337"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "'13"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.06996"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "'13"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02808"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("G12").Value = "'13"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09365"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("G13").Value = "'13"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001511"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("G14").Value = "'13"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005985"
$ws.Range("E15").Value = "14OneONEWorstin24h"
$ws.Range("G15").Value = "'13"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006160"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("G16").Value = "'13"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.500"
$ws.Range("E17").Value = "16LEOLEO"
$ws.Range("G17").Value = "'13"
$ws.Range("D18").Value = "'2.121"
$ws.Range("G18").Value = "'13"
$ws.Range("D19").Value = "'0.3197"
$ws.Range("G19").Value = "'13"
$ws.Range("D20").Value = "'0.03196"
$ws.Range("G20").Value = "'13"
$ws.Range("G21").Value = "'13"
$ws.Range("D22").Value = "'3.759"
$ws.Range("G22").Value = "'13"
$ws.Range("D23").Value = "'0.04661"
$ws.Range("G23").Value = "'13"
$ws.Range("G24").Value = "'13"
$ws.Range("D25").Value = "'0.001234"
$ws.Range("G25").Value = "'13"
$ws.Range("D26").Value = "'0.004258"
$ws.Range("G26").Value = "'13"
$ws.Range("D27").Value = "'0.00009691"
$ws.Range("G27").Value = "'13"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("G28").Value = "'13"
$ws.Range("G29").Value = "'13"
$ws.Range("G30").Value = "'13"
$ws.Range("G31").Value = "'13"
$ws.Range("G32").Value = "'13"
$ws.Range("G33").Value = "'13"
$ws.Range("G34").Value = "'13"
$ws.Range("G35").Value = "'13"
$ws.Range("G36").Value = "'13"
$ws.Range("G37").Value = "'13"
$ws.Range("G38").Value = "'13"
$ws.Range("G39").Value = "'13"
$ws.Range("D40").Value = "'0.03618"
$ws.Range("G40").Value = "'13"
$ws.Range("D41").Value = "'0.006334"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'13"
$ws.Range("D42").Value = "'0.1049"
$ws.Range("G42").Value = "'13"
$ws.Range("D43").Value = "'0.002997"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'13"
$ws.Range("D44").Value = "'0.007295"
$ws.Range("G44").Value = "'13"
$ws.Range("D45").Value = "'0.00005279"
$ws.Range("G45").Value = "'13"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("G46").Value = "'13"
$ws.Range("D47").Value = "'0.1898"
$ws.Range("G47").Value = "'13"
$ws.Range("D48").Value = "'0.002299"
$ws.Range("G48").Value = "'13"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("G49").Value = "'13"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("G50").Value = "'13"
$ws.Range("G51").Value = "'13"
